$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.139235
$ws.Range("H2").Value = 63.417705
$ws.Range("I2").Value = 0.1633331201667119
$ws.Range("J2").Value = 0.1633331201667119
$ws.Range("M2").Value = 1.947351333333333
$ws.Range("N2").Value = 5.842054
$ws.Range("O2").Value = 0.009881137208588643
$ws.Range("P2").Value = 0.009881137208588641
$ws.Range("Q2").Value = 41.16551746289667
$ws.Range("R2").Value = 370.48965716607
$ws.Range("S2").Value = 0.001613916971074177
$ws.Range("T2").Value = 0.001613916971074177
$ws.Range("G3").Value = 21.139235
$ws.Range("H3").Value = 63.417705
$ws.Range("I3").Value = 0.1633331201667119
$ws.Range("J3").Value = 0.1633331201667119
$ws.Range("O3").Value = 0.5240248869075229
$ws.Range("P3").Value = 0.5240248869075228
$ws.Range("Q3").Value = 2183.124793999825
$ws.Range("R3").Value = 19648.12314599842
$ws.Range("S3").Value = 0.08559061982361406
$ws.Range("T3").Value = 0.08559061982361406
$ws.Range("G4").Value = 21.139235
$ws.Range("H4").Value = 63.417705
$ws.Range("I4").Value = 0.1633331201667119
$ws.Range("J4").Value = 0.1633331201667119
$ws.Range("M4").Value = 31.14585733333333
$ws.Range("N4").Value = 93.43757199999999
$ws.Range("O4").Value = 0.1580385031308132
$ws.Range("P4").Value = 0.1580385031308132
$ws.Range("Q4").Value = 658.3995974458065
$ws.Range("R4").Value = 5925.596377012259
$ws.Range("S4").Value = 0.0258129218228324
$ws.Range("T4").Value = 0.0258129218228324
$ws.Range("G5").Value = 21.139235
$ws.Range("H5").Value = 63.417705
$ws.Range("I5").Value = 0.1633331201667119
$ws.Range("J5").Value = 0.1633331201667119
$ws.Range("M5").Value = 60.71084966666666
$ws.Range("N5").Value = 182.132549
$ws.Range("O5").Value = 0.3080554727530752
$ws.Range("P5").Value = 0.3080554727530752
$ws.Range("Q5").Value = 1283.380918153338
$ws.Range("R5").Value = 11550.42826338004
$ws.Range("S5").Value = 0.05031566154919127
$ws.Range("T5").Value = 0.05031566154919128
$ws.Range("I6").Value = 0.5748271090353965
$ws.Range("J6").Value = 0.5748271090353966
$ws.Range("M6").Value = 1.947351333333333
$ws.Range("N6").Value = 5.842054
$ws.Range("O6").Value = 0.009881137208588643
$ws.Range("P6").Value = 0.009881137208588641
$ws.Range("Q6").Value = 144.8760384359918
$ws.Range("R6").Value = 1303.884345923926
$ws.Range("S6").Value = 0.005679945535595097
$ws.Range("T6").Value = 0.005679945535595097
$ws.Range("I7").Value = 0.5748271090353965
$ws.Range("J7").Value = 0.5748271090353966
$ws.Range("O7").Value = 0.5240248869075229
$ws.Range("P7").Value = 0.5240248869075228
$ws.Range("S7").Value = 0.301223710803652
$ws.Range("T7").Value = 0.301223710803652
$ws.Range("I8").Value = 0.5748271090353965
$ws.Range("J8").Value = 0.5748271090353966
$ws.Range("M8").Value = 31.14585733333333
$ws.Range("N8").Value = 93.43757199999999
$ws.Range("O8").Value = 0.1580385031308132
$ws.Range("P8").Value = 0.1580385031308132
$ws.Range("Q8").Value = 2317.141415063563
$ws.Range("R8").Value = 20854.27273557207
$ws.Range("S8").Value = 0.09084481587096684
$ws.Range("T8").Value = 0.09084481587096685
$ws.Range("I9").Value = 0.5748271090353965
$ws.Range("J9").Value = 0.5748271090353966
$ws.Range("M9").Value = 60.71084966666666
$ws.Range("N9").Value = 182.132549
$ws.Range("O9").Value = 0.3080554727530752
$ws.Range("P9").Value = 0.3080554727530752
$ws.Range("Q9").Value = 4516.67207618573
$ws.Range("R9").Value = 40650.04868567158
$ws.Range("S9").Value = 0.1770786368251825
$ws.Range("T9").Value = 0.1770786368251826
$ws.Range("G10").Value = 4.054539666666667
$ws.Range("H10").Value = 12.163619
$ws.Range("I10").Value = 0.03132755819197652
$ws.Range("J10").Value = 0.03132755819197652
$ws.Range("M10").Value = 1.947351333333333
$ws.Range("N10").Value = 5.842054
$ws.Range("O10").Value = 0.009881137208588643
$ws.Range("P10").Value = 0.009881137208588641
$ws.Range("Q10").Value = 7.895613225936223
$ws.Range("R10").Value = 71.060519033426
$ws.Range("S10").Value = 0.0003095519009049652
$ws.Range("T10").Value = 0.0003095519009049651
$ws.Range("G11").Value = 4.054539666666667
$ws.Range("H11").Value = 12.163619
$ws.Range("I11").Value = 0.03132755819197652
$ws.Range("J11").Value = 0.03132755819197652
$ws.Range("O11").Value = 0.5240248869075229
$ws.Range("P11").Value = 0.5240248869075228
$ws.Range("Q11").Value = 418.7268874467683
$ws.Range("R11").Value = 3768.541987020915
$ws.Range("S11").Value = 0.01641642013863934
$ws.Range("T11").Value = 0.01641642013863934
$ws.Range("G12").Value = 4.054539666666667
$ws.Range("H12").Value = 12.163619
$ws.Range("I12").Value = 0.03132755819197652
$ws.Range("J12").Value = 0.03132755819197652
$ws.Range("M12").Value = 31.14585733333333
$ws.Range("N12").Value = 93.43757199999999
$ws.Range("O12").Value = 0.1580385031308132
$ws.Range("P12").Value = 0.1580385031308132
$ws.Range("Q12").Value = 126.2821140103409
$ws.Range("R12").Value = 1136.539026093068
$ws.Range("S12").Value = 0.004950960403403416
$ws.Range("T12").Value = 0.004950960403403416
$ws.Range("G13").Value = 4.054539666666667
$ws.Range("H13").Value = 12.163619
$ws.Range("I13").Value = 0.03132755819197652
$ws.Range("J13").Value = 0.03132755819197652
$ws.Range("M13").Value = 60.71084966666666
$ws.Range("N13").Value = 182.132549
$ws.Range("O13").Value = 0.3080554727530752
$ws.Range("P13").Value = 0.3080554727530752
$ws.Range("Q13").Value = 246.1545481705368
$ws.Range("R13").Value = 2215.390933534831
$ws.Range("S13").Value = 0.0096506257490288
$ws.Range("T13").Value = 0.0096506257490288
$ws.Range("G14").Value = 29.83382566666667
$ws.Range("H14").Value = 89.50147699999999
$ws.Range("I14").Value = 0.2305122126059151
$ws.Range("J14").Value = 0.2305122126059151
$ws.Range("M14").Value = 1.947351333333333
$ws.Range("N14").Value = 5.842054
$ws.Range("O14").Value = 0.009881137208588643
$ws.Range("P14").Value = 0.009881137208588641
$ws.Range("Q14").Value = 58.09694019041756
$ws.Range("R14").Value = 522.872461713758
$ws.Range("S14").Value = 0.002277722801014403
$ws.Range("T14").Value = 0.002277722801014403
$ws.Range("G15").Value = 29.83382566666667
$ws.Range("H15").Value = 89.50147699999999
$ws.Range("I15").Value = 0.2305122126059151
$ws.Range("J15").Value = 0.2305122126059151
$ws.Range("O15").Value = 0.5240248869075229
$ws.Range("P15").Value = 0.5240248869075228
$ws.Range("Q15").Value = 3081.046429199938
$ws.Range("R15").Value = 27729.41786279944
$ws.Range("S15").Value = 0.1207941361416175
$ws.Range("T15").Value = 0.1207941361416175
$ws.Range("G16").Value = 29.83382566666667
$ws.Range("H16").Value = 89.50147699999999
$ws.Range("I16").Value = 0.2305122126059151
$ws.Range("J16").Value = 0.2305122126059151
$ws.Range("M16").Value = 31.14585733333333
$ws.Range("N16").Value = 93.43757199999999
$ws.Range("O16").Value = 0.1580385031308132
$ws.Range("P16").Value = 0.1580385031308132
$ws.Range("Q16").Value = 929.2000779215381
$ws.Range("R16").Value = 8362.800701293843
$ws.Range("S16").Value = 0.0364298050336106
$ws.Range("T16").Value = 0.0364298050336106
$ws.Range("G17").Value = 29.83382566666667
$ws.Range("H17").Value = 89.50147699999999
$ws.Range("I17").Value = 0.2305122126059151
$ws.Range("J17").Value = 0.2305122126059151
$ws.Range("M17").Value = 60.71084966666666
$ws.Range("N17").Value = 182.132549
$ws.Range("O17").Value = 0.3080554727530752
$ws.Range("P17").Value = 0.3080554727530752
$ws.Range("Q17").Value = 1811.236905030541
$ws.Range("R17").Value = 16301.13214527487
$ws.Range("S17").Value = 0.07101054862967253
$ws.Range("T17").Value = 0.07101054862967254
